# BardBotDispatcher/Data/Config.xlsx
# "Have separate python bots for each automation"
#
# Inserts a new settings row (BardBotPath) directly beneath the existing
# OrchestratorQueueName row on the Settings sheet, pushing the remaining
# rows down by one, and widens the Value column so the new path is
# readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Shift row 3 downward (and everything below it) to make room for the
# new setting, right after the existing "OrchestratorQueueName" row.
$ws.Rows(3).Insert()

# Name / Description / Value for the new BardBotPath setting (filled in
# that order, matching how the shared-string table was populated).
$ws.Range("A3").Value = "BardBotPath"
$ws.Range("C3").Value = "Path to BardBot python bot."
$ws.Range("B3").Value = "D:\Revature\220425-UiPath\music-majors\BardBotDispatcher\DispatcherPythonBot\dist\BardBotD.exe"

# Keep the new row's height the same as its neighbours.
$ws.Rows(3).RowHeight = 14.25

# The Value column needs to be widened to fit the new (long) path -
# best-fit width for the longest value now in the column (~94.7 chars).
$ws.Columns(2).ColumnWidth = 93.83

Write-Output "Inserted BardBotPath setting row"
